$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.948.66'
$ws.Range('E2').Value = '  -1.25%  '
$ws.Range('D3').Value = '3.155.83'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '''589.04'
$ws.Range('E5').Value = '  -1.76%  '
$ws.Range('D6').Value = '''138.00'
$ws.Range('E6').Value = '  -3.22%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '3.144.85'
$ws.Range('E8').Value = '  +0.80%  '
$ws.Range('D9').Value = '''0.516'
$ws.Range('E9').Value = '  -1.13%  '
$ws.Range('E10').Value = '  -1.66%  '
$ws.Range('D11').Value = '''5.29'
$ws.Range('E11').Value = '  -1.36%  '
$ws.Range('D12').Value = '''0.457'
$ws.Range('E12').Value = '  -1.97%  '
$ws.Range('D13').Value = '''0.0000244'
$ws.Range('E13').Value = '  -3.68%  '
$ws.Range('D14').Value = '''34.02'
$ws.Range('E14').Value = '  -3.00%  '
$ws.Range('D15').Value = '3.676.47'
$ws.Range('E15').Value = '  +1.13%  '
$ws.Range('E16').Value = '  +0.86%  '
$ws.Range('D17').Value = '3.156.57'
$ws.Range('E17').Value = '  +0.66%  '
$ws.Range('D18').Value = '62.956.35'
$ws.Range('E18').Value = '  -1.26%  '
$ws.Range('D19').Value = '''6.64'
$ws.Range('E19').Value = '  -2.66%  '
$ws.Range('D20').Value = '''475.80'
$ws.Range('E20').Value = '  -1.46%  '
$ws.Range('D21').Value = '''13.94'
$ws.Range('E21').Value = '  -5.64%  '
$ws.Range('D22').Value = '''0.699'
$ws.Range('E22').Value = '  -1.14%  '
$ws.Range('D23').Value = '''7.70'
$ws.Range('E23').Value = '  +1.27%  '
$ws.Range('D24').Value = '''84.50'
$ws.Range('E24').Value = '  -2.41%  '
$ws.Range('D25').Value = '''12.94'
$ws.Range('E25').Value = '  -3.07%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').Value = '''2.70'
$ws.Range('E27').Value = '  -1.30%  '
$ws.Range('D28').Value = '''7.01'
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('D29').Value = '''7.90'
$ws.Range('E29').Value = '  -3.81%  '
$ws.Range('D30').Value = '''2.08'
$ws.Range('E30').Value = '  +1.66%  '
$ws.Range('E31').Value = '  +0.14%  '
$ws.Range('D32').Value = '''26.74'
$ws.Range('E32').Value = '  -1.08%  '
$ws.Range('E33').Value = '  -4.15%  '
$ws.Range('D34').Value = '''2.52'
$ws.Range('E34').Value = '  -5.48%  '
$ws.Range('D35').Value = '''1.07'
$ws.Range('E35').Value = '  -2.87%  '
$ws.Range('D36').Value = '''52.45'
$ws.Range('E36').Value = '  -0.30%  '
$ws.Range('D37').Value = '''5.76'
$ws.Range('E37').Value = '  -3.46%  '
$ws.Range('D38').Value = '0.0₃0698'
$ws.Range('E38').Value = '  -5.68%  '
$ws.Range('D39').Value = '''0.0386'
$ws.Range('E39').Value = '  -1.93%  '
$ws.Range('D40').Value = '''417.01'
$ws.Range('E40').Value = '  -4.61%  '
$ws.Range('D41').Value = '''2.75'
$ws.Range('E41').Value = '  -6.54%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.947.46'
$ws.Range('E42').Value = '  +2.81%  '
$ws.Range('B43').Value = 'Cosmos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D43').Value = '''8.25'
$ws.Range('E43').Value = '  -0.16%  '
$ws.Range('E44').Value = '  -7.49%  '
$ws.Range('D45').Value = '''0.258'
$ws.Range('E45').Value = '  -0.45%  '
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').Value = '''2.12'
$ws.Range('E47').Value = '  -3.68%  '
$ws.Range('D48').Value = '''25.36'
$ws.Range('E48').Value = '  -1.89%  '
$ws.Range('E49').Value = '  -0.57%  '
$ws.Range('D50').Value = '''2.23'
$ws.Range('E50').Value = '  -6.28%  '
$ws.Range('D51').Value = '''119.99'
$ws.Range('E51').Value = '  -1.54%  '
